# AFDP-7308 Combine Transcribe and OCR processing into a single media
# processing module.
#
# The "OCR Workflow Rules" rule table referenced a standalone OCR
# business-process model/workflow. Update it to reference the new,
# combined Media Engine business-process model/workflow instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rule table "global" declaration: $model: OCRBusinessProcessModel -> $model: MediaEngineBusinessProcessModel
$ws.Range("C14").Value = "`$model: MediaEngineBusinessProcessModel"

# Drools import: OCR model class -> Media Engine model class
$ws.Range("D3").Value = "com.armedia.acm.services.mediaengine.model.MediaEngineBusinessProcessModel"

# Process name started by the Automatic/Manual OCR rows: OCRWorkflow -> MediaEngineWorkFlow
$ws.Range("E17").Value = "MediaEngineWorkFlow"
$ws.Range("E18").Value = "MediaEngineWorkFlow"

# Refresh the sheet's selection/scroll position to match the saved view.
$ws.Activate()
$ws.Range("E20").Select()
